# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Update the Rule column value in row 11 ("R40" -> "1")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# B11 currently holds the text "R40"; change it to the text "1".
# A leading apostrophe keeps Excel from coercing the numeric-looking
# text into a real number, so the cell keeps storing a shared string.
$ws.Range("B11").Value = "'1"
